$d = $word.ActiveDocument

# Locate the target paragraph by its distinctive original text.
$targetStart = -1
$targetEnd = -1
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*remove the LAST_FM_API_KEY from the js file*") {
        $targetStart = $p.Range.Start
        $targetEnd = $p.Range.End
        break
    }
}

if ($targetStart -eq -1) {
    throw "target paragraph not found"
}

# Clear the paragraph's text content (leave the paragraph mark itself, i.e.
# stop one character short of $targetEnd) so we can rebuild it with the
# exact run / proofErr / symbol structure the edit calls for.
$clearRng = $d.Range($targetStart, $targetEnd - 1)
$clearRng.Text = ""

# Re-insert the paragraph body using WordprocessingML so the run
# boundaries, <w:proofErr/> spell/grammar markers and the Wingdings
# arrow <w:sym/> come out exactly as authored.
$insertRng = $d.Range($targetStart, $targetStart)
$paraXml = '<w:p w14:paraId="774432D2" w14:textId="0B692DF7" w:rsidR="00C27085" w:rsidRDefault="00BA03EF" w:rsidP="00C27085"><w:r><w:t>5</w:t></w:r><w:r w:rsidR="00C27085"><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidR="00C27085"><w:t>replace</w:t></w:r><w:r w:rsidR="00C27085"><w:t xml:space="preserve"> the</w:t></w:r><w:r w:rsidR="00C27085"><w:t xml:space="preserve"> actual key with</w:t></w:r><w:r w:rsidR="00C27085"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C27085" w:rsidRPr="00C27085"><w:t>LAST_FM_API_KEY</w:t></w:r><w:r w:rsidR="00C27085"><w:t xml:space="preserve"> from the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00C27085"><w:t>js</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00C27085"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00C27085"><w:t>file</w:t></w:r><w:r w:rsidR="00C27085"><w:t xml:space="preserve">  [</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00C27085"><w:t xml:space="preserve">should be like this </w:t></w:r><w:r w:rsidR="00C27085"><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r w:rsidR="00C27085"><w:t xml:space="preserve">   </w:t></w:r><w:r w:rsidR="00C27085"><w:t>${</w:t></w:r><w:r w:rsidR="00C27085"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C27085"><w:t>LAST_FM_API_KEY</w:t></w:r><w:r w:rsidR="00C27085"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C27085"><w:t>}</w:t></w:r><w:r w:rsidR="00C27085"><w:t xml:space="preserve">       ]</w:t></w:r><w:r w:rsidR="00C27085"><w:br/><w:t xml:space="preserve">so there is no </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00C27085"><w:t>api</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00C27085"><w:t xml:space="preserve"> key</w:t></w:r><w:r w:rsidR="00C27085"><w:t xml:space="preserve"> in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00C27085"><w:t>js</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00C27085"><w:t xml:space="preserve"> file. It just in the html file</w:t></w:r><w:r w:rsidR="00C27085"><w:br/></w:r></w:p>'
$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $paraXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRng.InsertXML($packageXml)

# Styles.xml: "Default Paragraph Font" should no longer be marked
# semiHidden (only unhideWhenUsed remains).
$style = $d.Styles("Default Paragraph Font")
$style.UnhideWhenUsed = $true

Write-Host "done"
